# Add "Sed" usage example (row 31) and "cmd ren" style for+sed entry (row 32)
# to the end of the command table on sheet 工作表1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 31: Sed / basic ops / delete-lines example (wrapped, taller row) ---
$ws.Range("A31").Value = "Sed"
$ws.Range("B31").Value = "basic ops"
$ws.Range("C31").Value = "Delete lines:" + [char]10 + "> sed '/^u/d' input.txt              //delete lines start with char 'u'"
$ws.Range("C31").WrapText = $true
$ws.Rows.Item(31).RowHeight = 31.5

# --- Row 32: for + sed / clear comments example ---
$ws.Range("A32").Value = "for + sed"
$ws.Range("B32").Value = "Clear all comment and save java source to a file with suffix 'Cleaned'"
$ws.Range("C32").Value = "for f in `$(find *.java -type f); do sed '/*/d' `$f > `${f}.Cleaned.java; done"

# Keep the selection where Excel leaves it after appending these rows.
$ws.Range("C30").Select()
